$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.08460788428783417, 0.9854296445846558, 0.03740128874778748, 0.9965887665748596),
    @(0.0220959484577179, 0.9977836608886719, 0.08866985142230988, 0.9739894270896912),
    @(0.01343870814889669, 0.9980504512786865, 0.01149542070925236, 0.9965887665748596),
    @(0.009787436574697495, 0.9979273080825806, 0.005867223720997572, 0.9988060593605042),
    @(0.00767605984583497, 0.9982146024703979, 0.004041307140141726, 0.9989766478538513),
    @(0.007612740155309439, 0.998153030872345, 0.003407995449379086, 0.9989766478538513),
    @(0.006386157125234604, 0.9982761740684509, 0.01221900805830956, 0.9974415898323059),
    @(0.006233484949916601, 0.9982351660728455, 0.002858529333025217, 0.9989766478538513),
    @(0.006530633196234703, 0.9981735944747925, 0.002585718641057611, 0.9990619421005249),
    @(0.006057605147361755, 0.998091459274292, 0.003212876617908478, 0.9988913536071777),
    @(0.005151795689016581, 0.9984608888626099, 0.002426594961434603, 0.9990619421005249),
    @(0.006611099932342768, 0.9979683756828308, 0.002884282730519772, 0.9990619421005249),
    @(0.005288407206535339, 0.9983172416687012, 0.00370407453738153, 0.9988913536071777),
    @(0.005093896761536598, 0.9982556700706482, 0.006957088597118855, 0.9990619421005249),
    @(0.006080134306102991, 0.998153030872345, 0.0027964161708951, 0.9990619421005249),
    @(0.005751732736825943, 0.9981735944747925, 0.002354457508772612, 0.9990619421005249),
    @(0.005125947296619415, 0.9982146024703979, 0.004322472959756851, 0.9988060593605042),
    @(0.005334658548235893, 0.9982966780662537, 0.001891778665594757, 0.9990619421005249),
    @(0.005699350032955408, 0.9982556700706482, 0.003095085266977549, 0.9988913536071777),
    @(0.004392130766063929, 0.9984198212623596, 0.00190135498996824, 0.9990619421005249),
    @(0.004470640327781439, 0.9984813928604126, 0.001983948284760118, 0.9991471767425537),
    @(0.005200617946684361, 0.9982146024703979, 0.00475058751180768, 0.9988913536071777),
    @(0.005193050019443035, 0.9982556700706482, 0.002951384056359529, 0.9988913536071777),
    @(0.004938268568366766, 0.9985018968582153, 0.002054528566077352, 0.9991471767425537),
    @(0.004541105590760708, 0.9983788132667542, 0.002350390655919909, 0.9996588826179504),
    @(0.004520933609455824, 0.9985429644584656, 0.002998686861246824, 0.999744176864624),
    @(0.00454885745421052, 0.9987276792526245, 0.002621467225253582, 0.9992324709892273),
    @(0.005501982755959034, 0.9985224604606628, 0.0072540701366961, 0.9991471767425537),
    @(0.005488427355885506, 0.9985634684562683, 0.002440281677991152, 0.9993177652359009),
    @(0.004104959778487682, 0.9988507628440857, 0.00535019813105464, 0.999744176864624),
    @(0.004681806545704603, 0.9986661076545715, 0.005765160080045462, 0.999744176864624),
    @(0.004253255669027567, 0.9987481832504272, 0.007325727492570877, 0.999744176864624),
    @(0.00492444634437561, 0.9985840320587158, 0.002373769646510482, 0.9998294115066528),
    @(0.004176911432296038, 0.9988097548484802, 0.004936868790537119, 0.9990619421005249),
    @(0.00499038165435195, 0.9986866116523743, 0.003082116600126028, 0.9999147057533264),
    @(0.004368199501186609, 0.99876868724823, 0.003653128165751696, 0.999744176864624),
    @(0.004816306289285421, 0.9987481832504272, 0.004005470313131809, 0.999744176864624),
    @(0.003927084617316723, 0.9989123344421387, 0.0178611408919096, 0.999744176864624),
    @(0.004363788291811943, 0.9989534020423889, 0.003531167982146144, 0.9998294115066528),
    @(0.004773812368512154, 0.9988097548484802, 0.004849810153245926, 0.9998294115066528),
    @(0.008088907226920128, 0.9986866116523743, 0.002972200280055404, 0.9999147057533264),
    @(0.00439781229943037, 0.9987276792526245, 0.002961742924526334, 0.9999147057533264),
    @(0.004253683611750603, 0.9989123344421387, 0.003003776771947742, 0.9999147057533264),
    @(0.004112687893211842, 0.9989534020423889, 0.003385127754881978, 0.9999147057533264),
    @(0.004885356873273849, 0.9987481832504272, 0.007544227875769138, 0.9996588826179504),
    @(0.005001910962164402, 0.9985840320587158, 0.00247120251879096, 0.9999147057533264),
    @(0.004486942198127508, 0.9988507628440857, 0.003325406461954117, 0.999744176864624),
    @(0.005073931068181992, 0.998645544052124, 0.002453385619446635, 0.9998294115066528),
    @(0.005058792885392904, 0.9986250400543213, 0.002011256758123636, 0.9998294115066528),
    @(0.004878455772995949, 0.9986045360565186, 0.002493765437975526, 0.9999147057533264)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt 4; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}
